$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Volume/Number and report-week date range (new weekly issue) ---
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# --- Row 14 ---
$ws.Range("D14").Copy($ws.Range("C14"))

# --- Row 15 ---
$ws.Range("L15").Copy($ws.Range("M15"))
$ws.Range("F15").Value = 1
$ws.Range("M15").Value = 300

# --- Row 16 ---
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = -23.809523809523
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = 14.285714285714
$ws.Range("N16").Value = -88.059701492537

# --- Row 17 ---
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -45.454545454545
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 27
$ws.Range("K17").Value = -14.814814814814
$ws.Range("L17").Value = -20.689655172413
$ws.Range("M17").Value = 91.666666666666
$ws.Range("N17").Value = -43.902439024390

# --- Row 18 ---
$ws.Range("G18").Copy($ws.Range("C18"))
$ws.Range("G18").Copy($ws.Range("D18"))
$ws.Range("L15").Copy($ws.Range("E18"))
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = -31.25
$ws.Range("L18").Value = -55.102040816326
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -91.760299625468

# --- Row 19 ---
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = -46.666666666666
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -6.382978723404
$ws.Range("I19").Value = 142
$ws.Range("J19").Value = 144
$ws.Range("K19").Value = -1.388888888888
$ws.Range("L19").Value = -1.388888888888
$ws.Range("M19").Value = -7.189542483660
$ws.Range("N19").Value = -66.190476190476

# --- Row 20 ---
$ws.Range("D14").Copy($ws.Range("C20"))
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -25
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = -36.363636363636
$ws.Range("N20").Value = -95.238095238095

# --- Row 21 ---
$ws.Range("C21").Value = 15
$ws.Range("E21").Value = -34.782608695652
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -8.108108108108
$ws.Range("I21").Value = 215
$ws.Range("J21").Value = 235
$ws.Range("K21").Value = -8.510638297872
$ws.Range("L21").Value = -10.788381742738
$ws.Range("M21").Value = -0.462962962962
$ws.Range("N21").Value = -78.754940711462

# --- Row 22 ---
$ws.Range("D14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 133.333333333333
$ws.Range("I22").Value = 15
$ws.Range("K22").Value = 87.5
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 36.363636363636

# --- Row 24 ---
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 26.666666666666
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 28.333333333333
$ws.Range("I24").Value = 225
$ws.Range("J24").Value = 207
$ws.Range("K24").Value = 8.695652173913
$ws.Range("L24").Value = -7.786885245901
$ws.Range("M24").Value = 49.006622516556

# --- Row 25 ---
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 8.333333333333
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 36.956521739130
$ws.Range("I25").Value = 184
$ws.Range("J25").Value = 162
$ws.Range("K25").Value = 13.580246913580
$ws.Range("L25").Value = -8.457711442786

# --- Row 26 ---
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 48
$ws.Range("J26").Value = 57
$ws.Range("K26").Value = -15.789473684210
$ws.Range("L26").Value = -5.882352941176
$ws.Range("M26").Value = -11.111111111111

# --- Row 27 ---
$ws.Range("F27").Value = 1

# --- Row 28 ---
$ws.Range("D14").Copy($ws.Range("C28"))
$ws.Range("E28").Value = -100
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 30
$ws.Range("L28").Value = -7.142857142857

# --- Row 31 ---
$ws.Range("L31").Value = 0
